# 自动更新Excel文件 - 2026-01-06 23:22:43
# Daily refresh: decrement the "剩余" (remaining) day-counter in column E for
# every data row. When the counter would drop to zero (or below) it means the
# current cycle has finished, so it is reset back to the row's "总天" (total
# days, column D) value and the "开始时间" (start date, column F) is rolled
# forward to the new reference day (2026-01-07). Rows whose start date is not
# a valid calendar date are left untouched (they fail the date check and are
# skipped), matching the source data which contains a malformed date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newToday = 20260107

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

function ConvertTo-Date($ymd) {
    # Returns a DateTime for a valid YYYYMMDD integer, or $null if invalid.
    if ($null -eq $ymd) { return $null }
    $s = [string]([int64]$ymd)
    if ($s.Length -ne 8) { return $null }
    $y = [int]$s.Substring(0, 4)
    $m = [int]$s.Substring(4, 2)
    $d = [int]$s.Substring(6, 2)
    try {
        return (Get-Date -Year $y -Month $m -Day $d)
    } catch {
        return $null
    }
}

for ($row = 2; $row -le $lastRow; $row++) {
    $totalDays = $ws.Cells.Item($row, 4).Value2
    $remaining = $ws.Cells.Item($row, 5).Value2
    $startDate = $ws.Cells.Item($row, 6).Value2

    if ($null -eq $remaining -or $null -eq $totalDays -or $null -eq $startDate) {
        continue
    }

    # Skip rows with an unparseable start date (data error) - leave as-is.
    $parsed = ConvertTo-Date $startDate
    if ($null -eq $parsed) {
        continue
    }

    $newRemaining = $remaining - 1
    if ($newRemaining -le 0) {
        $ws.Cells.Item($row, 5).Value = $totalDays
        $ws.Cells.Item($row, 6).Value = $newToday
    } else {
        $ws.Cells.Item($row, 5).Value = $newRemaining
    }
}
